# "Update HPBar and add crystaldata, leveldata"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DatasPath")

# --- Rename the NPCPoolsData path entry to NPCsData (row 10) ---
$ws.Range("A10").Value = "NPCsData"
$ws.Range("B10").Value = "NPCsData.xlsx"

# --- Add two new data-file rows: CrystalsData (row 11) and LevelData
#     (row 12). Copy row 10's formatting down first so the new rows match
#     the rest of the table, then fill in their values. ---
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B12").PasteSpecial(-4122)

$ws.Range("A11").Value = "CrystalsData"
$ws.Range("B11").Value = "CrystalsData.xlsx"
$ws.Range("A12").Value = "LevelData"
$ws.Range("B12").Value = "LevelData.xlsx"

# --- Update HPBar: highlight the header row with an orange fill
#     (RGB 255,153,0 -> OLE/VBA color = R + G*256 + B*65536 = 39423) ---
$ws.Range("A1:B1").Interior.Color = 39423
